$d = $word.ActiveDocument

# Paragraph 2: citation tag update (Ref-AB12CD) -> (Ref-f665654)
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Find.Execute("(Ref-AB12CD)", $true, $false, $false, $false, $false, $true, 1, $false, "(Ref-f665654)", 2)

# Paragraph 3: citation tag update (Ref-JHD73K) -> (Lee 208)
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Find.Execute("(Ref-JHD73K)", $true, $false, $false, $false, $false, $true, 1, $false, "(Lee 208)", 2)

# Paragraph 5: citation tag update (Ref-J7X2BZ) -> (Johnson 45)
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute("(Ref-J7X2BZ)", $true, $false, $false, $false, $false, $true, 1, $false, "(Johnson 45)", 2)

# Paragraph 6: citation tag update (Ref-J7Y3X2) -> (Ref-s241817)
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Find.Execute("(Ref-J7Y3X2)", $true, $false, $false, $false, $false, $true, 1, $false, "(Ref-s241817)", 2)

# Paragraph 7: citation tag update (Ref-JHD73K) -> (Ref-u881103)
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Find.Execute("(Ref-JHD73K)", $true, $false, $false, $false, $false, $true, 1, $false, "(Ref-u881103)", 2)
